$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.91"
$ws.Range("E2").Value = "'0.89%"
$ws.Range("G2").Value = "'8"
$ws.Range("D3").Value = "'35.97"
$ws.Range("E3").Value = "'-3.94%"
$ws.Range("G3").Value = "'8"
$ws.Range("D4").Value = "'5.090"
$ws.Range("E4").Value = "'1.63%"
$ws.Range("G4").Value = "'8"
$ws.Range("D5").Value = "'0.07859"
$ws.Range("E5").Value = "'0.45%"
$ws.Range("G5").Value = "'8"
$ws.Range("D6").Value = "'2.116"
$ws.Range("E6").Value = "'-3.22%"
$ws.Range("G6").Value = "'8"
$ws.Range("D7").Value = "'7.966"
$ws.Range("E7").Value = "'-0.77%"
$ws.Range("G7").Value = "'8"
$ws.Range("D8").Value = "'4.099"
$ws.Range("E8").Value = "'1.94%"
$ws.Range("G8").Value = "'8"
$ws.Range("D9").Value = "'0.9187"
$ws.Range("E9").Value = "'1.03%"
$ws.Range("G9").Value = "'8"
$ws.Range("D10").Value = "'0.09699"
$ws.Range("E10").Value = "'0.11%"
$ws.Range("G10").Value = "'8"
$ws.Range("D11").Value = "'0.1862"
$ws.Range("E11").Value = "'-1.81%"
$ws.Range("G11").Value = "'8"
$ws.Range("D12").Value = "'0.08645"
$ws.Range("E12").Value = "'1.84%"
$ws.Range("G12").Value = "'8"
$ws.Range("E13").Value = "'-0.77%"
$ws.Range("G13").Value = "'8"
$ws.Range("D14").Value = "'0.09945"
$ws.Range("E14").Value = "'-0.15%"
$ws.Range("G14").Value = "'8"
$ws.Range("D15").Value = "'0.001432"
$ws.Range("E15").Value = "'-3.32%"
$ws.Range("G15").Value = "'8"
$ws.Range("D16").Value = "'0.005722"
$ws.Range("E16").Value = "'0.79%"
$ws.Range("G16").Value = "'8"
$ws.Range("D17").Value = "'3.462"
$ws.Range("E17").Value = "'-0.04%"
$ws.Range("G17").Value = "'8"
$ws.Range("D18").Value = "'2.494"
$ws.Range("E18").Value = "'20.50%"
$ws.Range("G18").Value = "'8"
$ws.Range("D19").Value = "'0.3426"
$ws.Range("E19").Value = "'-1.06%"
$ws.Range("G19").Value = "'8"
$ws.Range("B20").Value = "'MCDex"
$ws.Range("C20").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'5.261"
$ws.Range("E20").Value = "'10.53%"
$ws.Range("G20").Value = "'8"
$ws.Range("B21").Value = "'ProBitToken"
$ws.Range("C21").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1301"
$ws.Range("E21").Value = "'0.55%"
$ws.Range("G21").Value = "'8"
$ws.Range("D22").Value = "'0.2202"
$ws.Range("E22").Value = "'-0.24%"
$ws.Range("G22").Value = "'8"
$ws.Range("D23").Value = "'0.04542"
$ws.Range("E23").Value = "'-1.87%"
$ws.Range("G23").Value = "'8"
$ws.Range("D24").Value = "'0.005085"
$ws.Range("E24").Value = "'5.91%"
$ws.Range("G24").Value = "'8"
$ws.Range("D25").Value = "'0.001235"
$ws.Range("E25").Value = "'0.35%"
$ws.Range("G25").Value = "'8"
$ws.Range("G26").Value = "'8"
$ws.Range("D27").Value = "'0.0004754"
$ws.Range("E27").Value = "'-0.07%"
$ws.Range("G27").Value = "'8"
$ws.Range("G28").Value = "'8"
$ws.Range("G29").Value = "'8"
$ws.Range("G30").Value = "'8"
$ws.Range("G31").Value = "'8"
$ws.Range("G32").Value = "'8"
$ws.Range("G33").Value = "'8"
$ws.Range("G34").Value = "'8"
$ws.Range("G35").Value = "'8"
$ws.Range("G36").Value = "'8"
$ws.Range("G37").Value = "'8"
$ws.Range("G38").Value = "'8"
$ws.Range("D39").Value = "'0.01851"
$ws.Range("E39").Value = "'5.62%"
$ws.Range("G39").Value = "'8"
$ws.Range("D40").Value = "'0.04740"
$ws.Range("E40").Value = "'0.47%"
$ws.Range("G40").Value = "'8"
$ws.Range("D41").Value = "'0.007531"
$ws.Range("E41").Value = "'-6.84%"
$ws.Range("G41").Value = "'8"
$ws.Range("D42").Value = "'0.1402"
$ws.Range("E42").Value = "'0.67%"
$ws.Range("G42").Value = "'8"
$ws.Range("D43").Value = "'0.007737"
$ws.Range("E43").Value = "'0.95%"
$ws.Range("G43").Value = "'8"
$ws.Range("D44").Value = "'0.002206"
$ws.Range("E44").Value = "'1.52%"
$ws.Range("G44").Value = "'8"
$ws.Range("D45").Value = "'0.01122"
$ws.Range("E45").Value = "'13.60%"
$ws.Range("G45").Value = "'8"
$ws.Range("D46").Value = "'0.00006378"
$ws.Range("E46").Value = "'4.77%"
$ws.Range("G46").Value = "'8"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("G47").Value = "'8"
$ws.Range("D48").Value = "'0.0005806"
$ws.Range("E48").Value = "'0.09%"
$ws.Range("G48").Value = "'8"
$ws.Range("D49").Value = "'41.54"
$ws.Range("E49").Value = "'379.09%"
$ws.Range("G49").Value = "'8"
$ws.Range("D50").Value = "'0.002002"
$ws.Range("E50").Value = "'-25.60%"
$ws.Range("G50").Value = "'8"
$ws.Range("D51").Value = "'0.00002102"
$ws.Range("E51").Value = "'-0.06%"
$ws.Range("G51").Value = "'8"
